$d = $word.ActiveDocument

# Find the index of the paragraph right after "Landing page básica..." --
# that paragraph (already empty in the original document) becomes the
# anchor we insert the new "url" / "Iguana Page" block after.
$anchorIdx = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx += 1
    if ($p.Range.Text -like "Landing page b*sica*") {
        $anchorIdx = $idx + 1
    }
}

$texts = @("url", "", "Iguana Page", "", "Página sobre iguanas hecha como ejercicio para estudio de posicionamiento en CSS con display block y elementos float.")

$insertIdx = $anchorIdx
foreach ($t in $texts) {
    $p = $d.Paragraphs($insertIdx)
    $p.Range.InsertParagraphAfter()
    $insertIdx += 1
    $newP = $d.Paragraphs($insertIdx)
    if ($t -ne "") {
        $newP.Range.Text = $t
    } else {
        # Work around the COM shim leaving a stray empty run behind when a
        # paragraph's Range.Text is never touched: type a placeholder
        # character and delete it so the paragraph ends up with no <w:r>
        # at all, matching a genuinely empty paragraph.
        $newP.Range.Text = "x"
        $delR = $d.Range($newP.Range.Start, $newP.Range.Start + 1)
        $delR.Delete()
    }
}
